# DemoFixture / ContactRowHandler sheet fix-up
#
# The "company" data originally lived bundled into the office/group cell
# (e.g. "ACME HQ", "ACME - 107 rue Saint Joan, 87654 Paris") with no
# dedicated company column. This edit:
#   - introduces a real "company" value ("ACME") in column C for the ACME
#     contacts/offices (rows 2-4, 9-11)
#   - renames the old combined office/address text in column A/B to short
#     office labels: "ACME HQ" -> "Global", the Paris address -> "Paris
#     Office", the Amiens address -> "Amiens Property"
#   - leaves the selected cell on B9, matching the saved workbook state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-4: "Management Board" / UK ACME contacts
$ws.Range("C2").Style = "Normal"
$ws.Range("C2").Value = "ACME"

$ws.Range("A2").Value = "Global"

$ws.Range("C3").Style = "Normal"
$ws.Range("C3").Value = "ACME"

$ws.Range("C4").Style = "Normal"
$ws.Range("C4").Value = "ACME"

# Row 9-10: Paris office contacts
$ws.Range("B9").Value = "Paris Office"

$ws.Range("C9").Style = "Normal"
$ws.Range("C9").Value = "ACME"

$ws.Range("C10").Style = "Normal"
$ws.Range("C10").Value = "ACME"

# Row 11: Amiens property contact
$ws.Range("B11").Value = "Amiens Property"

$ws.Range("C11").Style = "Normal"
$ws.Range("C11").Value = "ACME"

# Restore the saved cursor/selection position
$ws.Range("B9").Select()
